$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Range("A1").Value = "Volume"
